# Generate Report for Handoff
# Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps and stamp the "ht" (handoff) priority on the rows whose
# handback is now up to date.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows on each locale sheet whose handoff just completed (row 9 is left
# untouched - its handoff is still pending).
$rows = @(8, 10, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column G.
    $overview.Range("G$r").Value = "2016-08-16 20:19:09"

    # zh-cn sheet: Priority column E becomes "ht"; Latest Handoff
    # Datetime column H gets the new timestamp.
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-16 20:18:58"

    # de-de sheet: same treatment, but it shares the Overview timestamp.
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-16 20:19:09"
}
